$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1847826086956522
$ws.Range("C2").Value = 0.5797101449275363
$ws.Range("J2").Value = 0.02173913043478261
$ws.Range("O2").Value = 0.003623188405797101
$ws.Range("P2").Value = 0.1413043478260869
$ws.Range("S2").Value = 0.06884057971014493
$ws.Range("B3").Value = 0.01829268292682927
$ws.Range("C3").Value = 0.01829268292682927
$ws.Range("J3").Value = 0.07317073170731707
$ws.Range("P3").Value = 0.7317073170731707
$ws.Range("S3").Value = 0.1585365853658537
$ws.Range("J4").Value = 0.06122448979591837
$ws.Range("P4").Value = 0.5918367346938775
$ws.Range("S4").Value = 0.3469387755102041
$ws.Range("B6").Value = 0.08230452674897119
$ws.Range("D6").Value = 0.02469135802469136
$ws.Range("F6").Value = 0.09465020576131687
$ws.Range("J6").Value = 0.2386831275720165
$ws.Range("O6").Value = 0.02880658436213992
$ws.Range("Q6").Value = 0.1440329218106996
$ws.Range("R6").Value = 0.05761316872427984
$ws.Range("S6").Value = 0.3292181069958848
$ws.Range("B7").Value = 0.08823529411764706
$ws.Range("D7").Value = 0.0196078431372549
$ws.Range("F7").Value = 0.04901960784313725
$ws.Range("J7").Value = 0.107843137254902
$ws.Range("O7").Value = 0.04901960784313725
$ws.Range("Q7").Value = 0.2254901960784314
$ws.Range("R7").Value = 0.07843137254901961
$ws.Range("S7").Value = 0.3823529411764706
$ws.Range("B8").Value = 0.09389671361502347
$ws.Range("D8").Value = 0.009389671361502348
$ws.Range("F8").Value = 0.0539906103286385
$ws.Range("J8").Value = 0.1408450704225352
$ws.Range("O8").Value = 0.009389671361502348
$ws.Range("Q8").Value = 0.176056338028169
$ws.Range("R8").Value = 0.1314553990610329
$ws.Range("S8").Value = 0.3849765258215962
$ws.Range("B9").Value = 0.04597701149425287
$ws.Range("D9").Value = 0.04022988505747126
$ws.Range("F9").Value = 0.08620689655172414
$ws.Range("J9").Value = 0.132183908045977
$ws.Range("O9").Value = 0.01724137931034483
$ws.Range("Q9").Value = 0.2068965517241379
$ws.Range("R9").Value = 0.09770114942528736
$ws.Range("S9").Value = 0.3735632183908046
$ws.Range("B10").Value = 0.1009946442234124
$ws.Range("D10").Value = 0.02448355011476664
$ws.Range("E10").Value = 0.001530221882172915
$ws.Range("F10").Value = 0.07192042846212701
$ws.Range("J10").Value = 0.1346595256312165
$ws.Range("O10").Value = 0.02295332823259373
$ws.Range("Q10").Value = 0.2096403978576894
$ws.Range("R10").Value = 0.09028309104820199
$ws.Range("S10").Value = 0.3435348125478194
$ws.Range("F11").Value = 0.003184713375796179
$ws.Range("G11").Value = 0.1560509554140127
$ws.Range("J11").Value = 0.09554140127388536
$ws.Range("K11").Value = 0.2038216560509554
$ws.Range("L11").Value = 0.5318471337579618
$ws.Range("S11").Value = 0.009554140127388535
$ws.Range("G12").Value = 0.711864406779661
$ws.Range("J12").Value = 0.1977401129943503
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("S12").Value = 0.03954802259887006
$ws.Range("G13").Value = 0.74
$ws.Range("J13").Value = 0.16
$ws.Range("S13").Value = 0.1
$ws.Range("F15").Value = 0.03389830508474576
$ws.Range("H15").Value = 0.1271186440677966
$ws.Range("I15").Value = 0.07203389830508475
$ws.Range("J15").Value = 0.3771186440677966
$ws.Range("K15").Value = 0.08898305084745763
$ws.Range("M15").Value = 0.00423728813559322
$ws.Range("O15").Value = 0.06779661016949153
$ws.Range("S15").Value = 0.2288135593220339
$ws.Range("F16").Value = 0.03278688524590164
$ws.Range("H16").Value = 0.2295081967213115
$ws.Range("I16").Value = 0.07103825136612021
$ws.Range("J16").Value = 0.3114754098360656
$ws.Range("K16").Value = 0.1366120218579235
$ws.Range("M16").Value = 0.0273224043715847
$ws.Range("O16").Value = 0.04371584699453552
$ws.Range("S16").Value = 0.1475409836065574
$ws.Range("F17").Value = 0.01720430107526882
$ws.Range("H17").Value = 0.1505376344086022
$ws.Range("I17").Value = 0.05376344086021505
$ws.Range("J17").Value = 0.4559139784946237
$ws.Range("K17").Value = 0.1075268817204301
$ws.Range("M17").Value = 0.02795698924731183
$ws.Range("O17").Value = 0.07956989247311828
$ws.Range("S17").Value = 0.1075268817204301
$ws.Range("F18").Value = 0.03167420814479638
$ws.Range("H18").Value = 0.1719457013574661
$ws.Range("I18").Value = 0.09502262443438914
$ws.Range("J18").Value = 0.416289592760181
$ws.Range("K18").Value = 0.09502262443438914
$ws.Range("M18").Value = 0.02714932126696833
$ws.Range("O18").Value = 0.06334841628959276
$ws.Range("S18").Value = 0.09954751131221719
$ws.Range("F19").Value = 0.02049530315969257
$ws.Range("H19").Value = 0.2126387702818104
$ws.Range("I19").Value = 0.08198121263877028
$ws.Range("J19").Value = 0.3757472245943638
$ws.Range("K19").Value = 0.1084543125533732
$ws.Range("M19").Value = 0.02391118701964133
$ws.Range("N19").Value = 0.0008539709649871904
$ws.Range("O19").Value = 0.06575576430401367
$ws.Range("S19").Value = 0.1101622544833476
